$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "22.442.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.565.21"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.002"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.56"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.98"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.137"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07480"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.98"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.979"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.944"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.562.74"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06726"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.353"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.16"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.02"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.47%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "22.428.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.389"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.617"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.77"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.025"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.98"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.737.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.048"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +5.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.135"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.007"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.598"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08285"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02436"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.334"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2258"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.42%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.336"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.17"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6148"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.88"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.764"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5761"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.029"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.41"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.223"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07331"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.49%  "
